# Add the new "2022-Q4" quarter to the workbook:
#  1. Insert a new worksheet named "2022-Q4" right after "总计" (copying the
#     "2022-Q3" sheet so formatting/styles match the other quarter sheets).
#  2. Populate it with the new quarter's fund data.
#  3. Prepend a corresponding new row to the "总计" (summary) sheet, shifting
#     the existing rows down by one.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)
$wsQ3_2022 = $wb.Worksheets.Item("2022-Q3")

# --- 1. Create the new "2022-Q4" sheet right after "总计" -------------------
$wsQ3_2022.Copy($null, $wsTotal)
$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

# --- 2. Fill in the new quarter's data --------------------------------------
function Set-TextCell($ws, $addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# Row 2 - fund A
Set-TextCell $wsQ4 "C2" "华宝英国富时100指数A"
Set-TextCell $wsQ4 "D2" "0.14"
Set-TextCell $wsQ4 "E2" "94.75"
Set-TextCell $wsQ4 "F2" "5.28"
Set-TextCell $wsQ4 "G2" "0.0074"
$wsQ4.Range("H2").Value = 3

# Row 3 - fund C
Set-TextCell $wsQ4 "C3" "华宝英国富时100指数C"
Set-TextCell $wsQ4 "D3" "0.08"
Set-TextCell $wsQ4 "E3" "94.75"
Set-TextCell $wsQ4 "F3" "5.28"
Set-TextCell $wsQ4 "G3" "0.0042"
$wsQ4.Range("H3").Value = 3

# --- 3. Prepend a row to "总计" for the new quarter -------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Range("A2:D2").Insert()
# The inserted row doesn't inherit the right per-column formatting; copy it
# from the row right below (the old first data row) instead.
$wsTotal.Range("A3:D3").Copy()
$wsTotal.Range("A2:D2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.01
